# Atualizando os dados , Modificando o para_km para 1247954.666
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldParaKm = 1253164.5
$newParaKm = 1247954.666

for ($row = 2; $row -le 23; $row++) {
    $dCell = $ws.Cells.Item($row, 4)   # column D
    $eCell = $ws.Cells.Item($row, 5)   # column E

    $oldD = [double]$dCell.Value2
    $oldE = [double]$eCell.Value2

    $milhoProducao = $oldD * $oldParaKm
    $bovinoProducao = $oldE * $oldParaKm

    $dCell.Value2 = $milhoProducao / $newParaKm
    $eCell.Value2 = $bovinoProducao / $newParaKm
}
